$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial value that was bumped by one day
# (45189 -> 45190) for every data row (rows 2 through 506).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 506 }

$range = $ws.Range("C2:C$lastRow")
$range.Value = 45190
